# Apply BOM updates: swap out obsolete quick-connect / wire rows for the
# updated connector part numbers and add a new "Primary side wire" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the new "Primary side wire" row + extra blank rows
#    between the Harness table and the Misc table (old row 22 -> row 26).
#    Insert 4 rows at row 20, then fully clear them so they don't carry
#    any stray formatting (Excel normally copies the format of the row
#    above when inserting).
# ---------------------------------------------------------------------
$ws.Range("A20:A23").EntireRow.Insert()
$ws.Range("A20:F23").Clear()

# ---------------------------------------------------------------------
# 2) Row 11 - "Quick Connect - Female (power input)"
#    now references the 3-520276-2 / A99878CT-ND connector.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Quick Connect - Female (power input)"
$ws.Range("B11").Value = "3-520276-2"
$ws.Range("C11").Value = "A99878CT-ND"
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = "0.187"" (4.75mm) Quick Connect Female 14-16 AWG Crimp Connector Fully Insulated"
$ws.Range("F11").Value = "0.8 mm thick"
$ws.Rows.Item(11).RowHeight = 17

# Give C11 the same Arial-13 note style used elsewhere (e.g. C5 / C28),
# reusing the existing cell format via copy/paste-special so no new
# style record is created.
$ws.Range("C5").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Hyperlink B11 to the new connector's Digi-Key page.
$ws.Range("B11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.digikey.ca/en/products/detail/te-connectivity-amp-connectors/3-520276-2/2060928") | Out-Null
$ws.Range("B11").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3) Row 12 - "Quick Connect - Female (power distribution)"
#    now references the 4-520448-2 / A27831-ND connector.
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "Quick Connect - Female (power distribution)"
$ws.Range("B12").Value = "4-520448-2"
$ws.Range("C12").Value = "A27831-ND"
$ws.Range("E12").Value = "0.250"" (6.35mm) Quick Connect Female 10-12 AWG Crimp Connector Fully Insulated"
$ws.Range("F12").Value = "0.8mm"

# ---------------------------------------------------------------------
# 4) Row 14 - "Quick Connect - Male (power)" gains Notes + thickness.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Quick Connect - Male (power)"
$ws.Range("B14").Value = "1217861-1"
$ws.Range("C14").Value = "A100452CT-ND"
$ws.Range("E14").Value = "0.250"" (6.35mm) Quick Connect Male  Solder Connector Non-Insulated"
$ws.Range("F14").Value = "0.8mm"

# ---------------------------------------------------------------------
# 5) Row 15 - "Power Distribution wire" becomes a bold sub-heading.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Power Distribution wire"
$ws.Range("A15").Font.Bold = $true

# Rows 16-18 (Black / White / Brown wire) are unchanged in content.

# ---------------------------------------------------------------------
# 6) Row 20 (new) - "Primary side wire"
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Primary side wire"
$ws.Range("B20").Value = "55A0111-14-9"
$ws.Range("C20").Value = "A132382-DS-ND"
$ws.Range("E20").Value = "14 AWG Hook-Up Wire 19/27 White 600V Enter Number of Feet in Order Quantity"

# ---------------------------------------------------------------------
# 7) Update the selected cell to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("E22").Select()
